# Update the NATMI LR-pair TPM output (Jam2-F11r) sheet with recomputed
# stats after the underlying TPM values changed.
# Columns E:T for data rows 2..17 are replaced with the refreshed values
# (ligand / receptor expression stats and derived edge-weight metrics).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object 'object[,]' 16,16
$data[0,0] = 3
$data[0,1] = 1
$data[0,2] = 115.5575153333333
$data[0,3] = 346.672546
$data[0,4] = 0.9048104954928987
$data[0,5] = 0.9048104954928987
$data[0,6] = 3
$data[0,7] = 1
$data[0,8] = 40.91514966666667
$data[0,9] = 122.745449
$data[0,10] = 0.8529192913871414
$data[0,11] = 0.8529192913871415
$data[0,12] = 4728.053034971463
$data[0,13] = 42552.47731474316
$data[0,14] = 0.7717303266554515
$data[0,15] = 0.7717303266554516
$data[1,0] = 3
$data[1,1] = 1
$data[1,2] = 115.5575153333333
$data[1,3] = 346.672546
$data[1,4] = 0.9048104954928987
$data[1,5] = 0.9048104954928987
$data[1,6] = 2
$data[1,7] = 0.6666666666666666
$data[1,8] = 0.165216
$data[1,9] = 0.495648
$data[1,10] = 0.00344410114086962
$data[1,11] = 0.003444101140869621
$data[1,12] = 19.091950453312
$data[1,13] = 171.827554079808
$data[1,14] = 0.003116258859797899
$data[1,15] = 0.0031162588597979
$data[2,0] = 3
$data[2,1] = 1
$data[2,2] = 115.5575153333333
$data[2,3] = 346.672546
$data[2,4] = 0.9048104954928987
$data[2,5] = 0.9048104954928987
$data[2,6] = 2
$data[2,7] = 0.6666666666666666
$data[2,8] = 0.4441646666666667
$data[2,9] = 1.332494
$data[2,10] = 0.009259079236881667
$data[2,11] = 0.009259079236881667
$data[2,12] = 51.32656527885823
$data[2,13] = 461.9390875097241
$data[2,14] = 0.008377712072130911
$data[2,15] = 0.008377712072130911
$data[3,0] = 3
$data[3,1] = 1
$data[3,2] = 115.5575153333333
$data[3,3] = 346.672546
$data[3,4] = 0.9048104954928987
$data[3,5] = 0.9048104954928987
$data[3,6] = 3
$data[3,7] = 1
$data[3,8] = 6.446186333333333
$data[3,9] = 19.338559
$data[3,10] = 0.1343775282351073
$data[3,11] = 0.1343775282351073
$data[3,12] = 744.9052760556905
$data[3,13] = 6704.147484501214
$data[3,14] = 0.1215861979055184
$data[3,15] = 0.1215861979055184
$data[4,0] = 3
$data[4,1] = 1
$data[4,2] = 5.519651666666666
$data[4,3] = 16.558955
$data[4,4] = 0.04321864090845719
$data[4,5] = 0.04321864090845719
$data[4,6] = 3
$data[4,7] = 1
$data[4,8] = 40.91514966666667
$data[4,9] = 122.745449
$data[4,10] = 0.8529192913871414
$data[4,11] = 0.8529192913871415
$data[4,12] = 225.8373740495328
$data[4,13] = 2032.536366445795
$data[4,14] = 0.03686201257835663
$data[4,15] = 0.03686201257835663
$data[5,0] = 3
$data[5,1] = 1
$data[5,2] = 5.519651666666666
$data[5,3] = 16.558955
$data[5,4] = 0.04321864090845719
$data[5,5] = 0.04321864090845719
$data[5,6] = 2
$data[5,7] = 0.6666666666666666
$data[5,8] = 0.165216
$data[5,9] = 0.495648
$data[5,10] = 0.00344410114086962
$data[5,11] = 0.003444101140869621
$data[5,12] = 0.9119347697599998
$data[5,13] = 8.20741292784
$data[5,14] = 0.0001488493704596519
$data[5,15] = 0.0001488493704596519
$data[6,0] = 3
$data[6,1] = 1
$data[6,2] = 5.519651666666666
$data[6,3] = 16.558955
$data[6,4] = 0.04321864090845719
$data[6,5] = 0.04321864090845719
$data[6,6] = 2
$data[6,7] = 0.6666666666666666
$data[6,8] = 0.4441646666666667
$data[6,9] = 1.332494
$data[6,10] = 0.009259079236881667
$data[6,11] = 0.009259079236881667
$data[6,12] = 2.451634242641111
$data[6,13] = 22.06470818377
$data[6,14] = 0.0004001648206817406
$data[6,15] = 0.0004001648206817406
$data[7,0] = 3
$data[7,1] = 1
$data[7,2] = 5.519651666666666
$data[7,3] = 16.558955
$data[7,4] = 0.04321864090845719
$data[7,5] = 0.04321864090845719
$data[7,6] = 3
$data[7,7] = 1
$data[7,8] = 6.446186333333333
$data[7,9] = 19.338559
$data[7,10] = 0.1343775282351073
$data[7,11] = 0.1343775282351073
$data[7,12] = 35.58070313842722
$data[7,13] = 320.226328245845
$data[7,14] = 0.00580761413895917
$data[7,15] = 0.005807614138959171
$data[8,0] = 3
$data[8,1] = 1
$data[8,2] = 6.580297333333334
$data[8,3] = 19.740892
$data[8,4] = 0.05152345196666309
$data[8,5] = 0.05152345196666309
$data[8,6] = 3
$data[8,7] = 1
$data[8,8] = 40.91514966666667
$data[8,9] = 122.745449
$data[8,10] = 0.8529192913871414
$data[8,11] = 0.8529192913871415
$data[8,12] = 269.233850244501
$data[8,13] = 2423.104652200509
$data[8,14] = 0.0439453461412257
$data[8,15] = 0.0439453461412257
$data[9,0] = 3
$data[9,1] = 1
$data[9,2] = 6.580297333333334
$data[9,3] = 19.740892
$data[9,4] = 0.05152345196666309
$data[9,5] = 0.05152345196666309
$data[9,6] = 2
$data[9,7] = 0.6666666666666666
$data[9,8] = 0.165216
$data[9,9] = 0.495648
$data[9,10] = 0.00344410114086962
$data[9,11] = 0.003444101140869621
$data[9,12] = 1.087170404224
$data[9,13] = 9.784533638016002
$data[9,14] = 0.0001774519796999254
$data[9,15] = 0.0001774519796999255
$data[10,0] = 3
$data[10,1] = 1
$data[10,2] = 6.580297333333334
$data[10,3] = 19.740892
$data[10,4] = 0.05152345196666309
$data[10,5] = 0.05152345196666309
$data[10,6] = 2
$data[10,7] = 0.6666666666666666
$data[10,8] = 0.4441646666666667
$data[10,9] = 1.332494
$data[10,10] = 0.009259079236881667
$data[10,11] = 0.009259079236881667
$data[10,12] = 2.922735571627556
$data[10,13] = 26.30462014464801
$data[10,14] = 0.0004770597243170001
$data[10,15] = 0.0004770597243170001
$data[11,0] = 3
$data[11,1] = 1
$data[11,2] = 6.580297333333334
$data[11,3] = 19.740892
$data[11,4] = 0.05152345196666309
$data[11,5] = 0.05152345196666309
$data[11,6] = 3
$data[11,7] = 1
$data[11,8] = 6.446186333333333
$data[11,9] = 19.338559
$data[11,10] = 0.1343775282351073
$data[11,11] = 0.1343775282351073
$data[11,12] = 42.41782273940311
$data[11,13] = 381.760404654628
$data[11,14] = 0.006923594121420463
$data[11,15] = 0.006923594121420465
$data[12,0] = 1
$data[12,1] = 0.3333333333333333
$data[12,2] = 0.057141
$data[12,3] = 0.171423
$data[12,4] = 0.0004474116319810314
$data[12,5] = 0.0004474116319810314
$data[12,6] = 3
$data[12,7] = 1
$data[12,8] = 40.91514966666667
$data[12,9] = 122.745449
$data[12,10] = 0.8529192913871414
$data[12,11] = 0.8529192913871415
$data[12,12] = 2.337932567103
$data[12,13] = 21.041393103927
$data[12,14] = 0.0003816060121076258
$data[12,15] = 0.0003816060121076258
$data[13,0] = 1
$data[13,1] = 0.3333333333333333
$data[13,2] = 0.057141
$data[13,3] = 0.171423
$data[13,4] = 0.0004474116319810314
$data[13,5] = 0.0004474116319810314
$data[13,6] = 2
$data[13,7] = 0.6666666666666666
$data[13,8] = 0.165216
$data[13,9] = 0.495648
$data[13,10] = 0.00344410114086962
$data[13,11] = 0.003444101140869621
$data[13,12] = 0.009440607456
$data[13,13] = 0.08496546710400001
$data[13,14] = 0.000001540930912144209
$data[13,15] = 0.000001540930912144209
$data[14,0] = 1
$data[14,1] = 0.3333333333333333
$data[14,2] = 0.057141
$data[14,3] = 0.171423
$data[14,4] = 0.0004474116319810314
$data[14,5] = 0.0004474116319810314
$data[14,6] = 2
$data[14,7] = 0.6666666666666666
$data[14,8] = 0.4441646666666667
$data[14,9] = 1.332494
$data[14,10] = 0.009259079236881667
$data[14,11] = 0.009259079236881667
$data[14,12] = 0.025380013218
$data[14,13] = 0.228420118962
$data[14,14] = 0.000004142619752014909
$data[14,15] = 0.000004142619752014909
$data[15,0] = 1
$data[15,1] = 0.3333333333333333
$data[15,2] = 0.057141
$data[15,3] = 0.171423
$data[15,4] = 0.0004474116319810314
$data[15,5] = 0.0004474116319810314
$data[15,6] = 3
$data[15,7] = 1
$data[15,8] = 6.446186333333333
$data[15,9] = 19.338559
$data[15,10] = 0.1343775282351073
$data[15,11] = 0.1343775282351073
$data[15,12] = 0.368341533273
$data[15,13] = 3.315073799457
$data[15,14] = 0.00006012206920924647
$data[15,15] = 0.00006012206920924649
$ws.Range("E2:T17").Value = $data
